# Revert "adding term 2.0.0"
# - Restore Version/Date/Contact values on the Metadata sheet
# - Remove the 12 concept rows that were inserted on the "Include from FSIII" sheet

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$wsMeta.Range("B10").Value = "No display for ContactDetail"

$wsInclude = $wb.Worksheets.Item("Include from FSIII")
$wsInclude.Range("2:13").Delete()
